# Updated jxls 3 report in demo
#
# The sample jxls3 template is re-pointed from the old "${results.*}"
# placeholders (driven directly off the "results" collection) to the
# jxls "jx:each" loop-variable style ("${row.*}"), and the jx:area /
# jx:each directives that drive the report are documented as cell
# comments on A1 and A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Data cells: ${results.ITEM_NAME} / ${results.VOLUME} -> ${row.*} ---
$ws.Range("A5").Value = '${row.ITEM_NAME}'
$ws.Range("B5").Value = '${row.VOLUME}'

# --- Selection shown in the saved view: A5:B5, active cell B5 ---
$ws.Range("A5:B5").Select()
$ws.Range("B5").Activate()

# --- jxls directive comments ---
$jxArea = $ws.Range("A1").AddComment("Author:" + [char]10 + 'jx:area(lastCell="B5")')
$jxArea.Shape.TextFrame.Characters(1,7).Font.Bold = $true
$jxArea.Visible = $false

$jxEach = $ws.Range("A5").AddComment("Author:" + [char]10 + 'jx:each(items="results" var="row" lastCell="B5")')
$jxEach.Shape.TextFrame.Characters(1,7).Font.Bold = $true
$jxEach.Visible = $false

# --- Refresh the chart so its cached category labels follow the renamed
#     placeholder text ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.Refresh()
